$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Login")
$ws2 = $wb.Worksheets.Item("Registration")
$ws3 = $wb.Worksheets.Item("AddUser")

# ---------------------------------------------------------------
# AddUser (sheet3) - populate the new test-data table
# ---------------------------------------------------------------
$ws3.Range("A1").Value = "tcid"
$ws3.Range("B1").Value = "tcinfo"
$ws3.Range("C1").Value = "Username"
$ws3.Range("D1").Value = "Mobile"
$ws3.Range("E1").Value = "Email"
$ws3.Range("F1").Value = "Cources"
$ws3.Range("G1").Value = "Gender"
$ws3.Range("H1").Value = "selectvalue"
$ws3.Range("I1").Value = "Password"

$ws3.Range("A2").Value = 3
$ws3.Range("B2").Value = "validInfo"
$ws3.Range("C2").Value = "subhash"

$ws3.Range("D2").NumberFormat = "0;[Red]0"
$ws3.Range("D2").Value = "8654852352"

$ws3.Range("E2").Value = "ingale.subhash@gmail.com"
$ws3.Range("F2").Value = "selenium"
$ws3.Range("G2").Value = "Male"
$ws3.Range("H2").Value = "HP"

$ws3.Range("I2").NumberFormat = "@"
$ws3.Range("I2").Value = "123456"

# Row 2 custom height (12.75pt, as authored by Excel for this row)
$ws3.Rows.Item(2).RowHeight = 12.75

# Column widths (best-fit values from the authored workbook)
$ws3.Columns.Item(1).ColumnWidth = 3.4518229166666665
$ws3.Columns.Item(2).ColumnWidth = 14.166666666666666
$ws3.Columns.Item(3).ColumnWidth = 9.166666666666666
$ws3.Columns.Item(4).ColumnWidth = 10.166666666666666
$ws3.Columns.Item(5).ColumnWidth = 25.022135416666668
$ws3.Columns.Item(6).ColumnWidth = 8.451822916666666
$ws3.Columns.Item(7).ColumnWidth = 6.736979166666667
$ws3.Columns.Item(8).ColumnWidth = 10.307291666666666
$ws3.Columns.Item(9).ColumnWidth = 8.592447916666666

# Borders around the populated table, matching the other sheets' "bordered cell" look
$ws1.Range("C3").Copy()
$ws3.Range("A1:C1").PasteSpecial(-4122)
$ws3.Range("H1:I1").PasteSpecial(-4122)
$ws3.Range("A2:C2").PasteSpecial(-4122)
$ws3.Range("F2:H2").PasteSpecial(-4122)

$ws1.Range("D3").Copy()
$ws3.Range("D1:G1").PasteSpecial(-4122)

# D2: border + quotePrefix + custom red-negative number format (style index 5 upstream)
$ws1.Range("D3").Copy()
$ws3.Range("D2").PasteSpecial(-4122)
$ws3.Range("D2").NumberFormat = "0;[Red]0"

# I2: border + quotePrefix (matches the existing "123456" style used elsewhere)
$ws1.Range("D3").Copy()
$ws3.Range("I2").PasteSpecial(-4122)

# E2: hyperlink-styled cell (font + border), then wired up as a real hyperlink
$ws1.Range("C4").Copy()
$ws3.Range("E2").PasteSpecial(-4122)
$ws3.Hyperlinks.Add($ws3.Range("E2"), "mailto:ingale.subhash@gmail.com")

# Re-apply the cell text/number-format choices that PasteSpecial(Formats) does not disturb,
# in case paste order above changed number formatting on D2/I2.
$ws3.Range("D2").NumberFormat = "0;[Red]0"
$ws3.Range("I2").NumberFormat = "@"

# Page setup + print orientation (adds <pageSetup orientation="portrait" .../>)
$ws3.PageSetup.Orientation = 1

# ---------------------------------------------------------------
# Selections / active sheet / active cell bookkeeping
# ---------------------------------------------------------------
$ws2.Range("G2").Select()

$ws3.Activate()
$ws3.Range("B10").Select()
